# Adjust Investment Summary table column widths for better formatting
# (also touches the Timeline & Milestones table on slide 5, whose grid
# widths/extent shift by a rounding-sized amount as a side effect of the
# same column-width pass).

$p = $ppt.ActivePresentation

# EMU -> point conversion factor used by the PowerPoint object model.
$EMU_PER_POINT = 12700

# --- Slide 5: "Timeline & Milestones" table --------------------------------
$slide5 = $p.Slides.Item(5)
$tbl5 = $slide5.Shapes.Item(3).Table

$tbl5.Columns.Item(1).Width = 871093 / $EMU_PER_POINT
$tbl5.Columns.Item(2).Width = 2177733 / $EMU_PER_POINT
$tbl5.Columns.Item(3).Width = 1306639 / $EMU_PER_POINT
$tbl5.Columns.Item(4).Width = 4355466 / $EMU_PER_POINT

# --- Slide 8: "Investment Summary" table ------------------------------------
$slide8 = $p.Slides.Item(8)
$tbl8 = $slide8.Shapes.Item(3).Table

$tbl8.Columns.Item(1).Width = 1742186 / $EMU_PER_POINT
$tbl8.Columns.Item(2).Width = 1045311 / $EMU_PER_POINT
$tbl8.Columns.Item(3).Width = 2003514 / $EMU_PER_POINT
$tbl8.Columns.Item(4).Width = 1132421 / $EMU_PER_POINT
$tbl8.Columns.Item(5).Width = 871093 / $EMU_PER_POINT
$tbl8.Columns.Item(6).Width = 871093 / $EMU_PER_POINT
$tbl8.Columns.Item(7).Width = 1045311 / $EMU_PER_POINT
